$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix A10's value (tiny correction in fractional time component)
$ws.Range("A10").Value = 45866.45858155093

# Append new row 11 with the new sensor reading
$ws.Range("A11").Value = 45866.50026385241
$ws.Range("B11").Value = 2025
$ws.Range("C11").Value = 31
$ws.Range("D11").Value = 18.33
$ws.Range("E11").Value = 77.43000000000001
$ws.Range("F11").Value = 471.09
$ws.Range("G11").Value = 12.94
$ws.Range("H11").Value = "ESE"
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = "12:00:22"

# Match the date-time number format used in column A for the other rows
$ws.Range("A11").NumberFormat = $ws.Range("A10").NumberFormat
